$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fundamentals_Remaining")

$rowsToZero = @(12,15,16,61,118,170,171,172,180,181,182,184,185,186,189,190,191,192,193,194,202,206,207,212,213,214,215,238,240,241,246,268,278,281,283,287,296,299,370,458,479,481,487,495,508,519,520,528,572,583,599,600,604,609,633,634,639,684,693,737,738,741,742,743,744,747,748,755,756,758,759,760,761,903,904,905,907,908,909,910)

foreach ($r in $rowsToZero) {
    $ws.Cells.Item($r, 4).Value = 0
}

$ws.Cells.Item(269, 4).Value = 1
